$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.916.77"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.088.59"
$ws.Range("E3").Value = "  +4.98%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.03"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.88"
$ws.Range("E6").Value = "  +6.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.084.91"
$ws.Range("E8").Value = "  +4.97%  "
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("E11").Value = "  +3.70%  "
$ws.Range("E12").Value = "  +5.43%  "
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.52"
$ws.Range("E14").Value = "  +6.50%  "
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.598.74"
$ws.Range("E16").Value = "  +4.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.919.93"
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.20"
$ws.Range("E18").Value = "  +3.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.086.60"
$ws.Range("E19").Value = "  +4.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.19"
$ws.Range("E20").Value = "  +8.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.43"
$ws.Range("E21").Value = "  +4.72%  "
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.37"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +7.47%  "
$ws.Range("E26").Value = "  +7.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.23"
$ws.Range("E33").Value = "  +4.11%  "
$ws.Range("E34").Value = "  +3.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("B38").Value = "Arweave"
$ws.Range("C38").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.23"
$ws.Range("E38").Value = "  +7.12%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.12"
$ws.Range("E39").Value = "  +7.16%  "
$ws.Range("E40").Value = "  +7.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.20"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.69"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "385.01"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.772.82"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.02"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.87"
$ws.Range("E50").Value = "  +7.14%  "
$ws.Range("E51").Value = "  +2.47%  "
